$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updated values (B2:E2, G2; F2 stays 1)
$ws.Range("B2").Value = 0.003078177322033415
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 32.88861328645544

# Row 3 updated values (B3:E3, G3; F3 stays 1)
$ws.Range("B3").Value = 1.445647641019636
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 4.327115817150455
